# fall 24 final inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 10.7

$ws.Range("D3").Value = 10.17
$ws.Range("E3").Value = 10.85

$ws.Range("C4").Value = 9.83
$ws.Range("F4").Value = 9.630000000000001

$ws.Range("B5").Value = 9.300000000000001
$ws.Range("C5").Value = 9.109999999999999
$ws.Range("F5").Value = 10.17

$ws.Range("D6").Value = 10.37
$ws.Range("E6").Value = 9.83
$ws.Range("G6").Value = 10.26
$ws.Range("H6").Value = 10.56

$ws.Range("F7").Value = 9.74

$ws.Range("F8").Value = 9.44
$ws.Range("J8").Value = 10.81

$ws.Range("H10").Value = 9.19
